# Updated charts: re-indexed the "Unnamed: 0" sort-order column (D) so that
# unidentified hits ("! Not Identified to Species") render as a single data
# slice in the pie/bar charts, and species within each LocSeason group are
# renumbered accordingly.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 168
$ws.Range("D3").Value = 8
$ws.Range("D4").Value = 9
$ws.Range("D6").Value = 10
$ws.Range("D7").Value = 11
$ws.Range("D8").Value = 12
$ws.Range("D9").Value = 2
$ws.Range("D10").Value = 4
$ws.Range("D11").Value = 5
$ws.Range("D12").Value = 1
$ws.Range("D13").Value = 6
$ws.Range("D15").Value = 92
$ws.Range("D17").Value = 93
$ws.Range("D18").Value = 82
$ws.Range("D19").Value = 85
$ws.Range("D20").Value = 83
$ws.Range("D21").Value = 94
$ws.Range("D22").Value = 95
$ws.Range("D23").Value = 90
$ws.Range("D24").Value = 87
$ws.Range("D25").Value = 91
$ws.Range("D26").Value = 80
$ws.Range("D28").Value = 81
$ws.Range("D29").Value = 88
$ws.Range("D31").Value = 89
$ws.Range("D33").Value = 150
$ws.Range("D35").Value = 151
$ws.Range("D37").Value = 146
$ws.Range("D38").Value = 143
$ws.Range("D39").Value = 142
$ws.Range("D40").Value = 152
$ws.Range("D41").Value = 147
$ws.Range("D42").Value = 148
$ws.Range("D43").Value = 149
$ws.Range("D46").Value = 38
$ws.Range("D48").Value = 39
$ws.Range("D49").Value = 35
$ws.Range("D50").Value = 32
$ws.Range("D52").Value = 40
$ws.Range("D53").Value = 41
$ws.Range("D54").Value = 33
$ws.Range("D55").Value = 36
$ws.Range("D56").Value = 31
$ws.Range("D57").Value = 37
$ws.Range("D58").Value = 29
$ws.Range("D62").Value = 108
$ws.Range("D63").Value = 109
$ws.Range("D66").Value = 104
$ws.Range("D67").Value = 106
$ws.Range("D68").Value = 102
$ws.Range("D69").Value = 107
$ws.Range("D70").Value = 103
$ws.Range("D71").Value = 105
$ws.Range("D74").Value = 168
$ws.Range("D76").Value = 169
$ws.Range("D77").Value = 170
$ws.Range("D80").Value = 171
$ws.Range("D81").Value = 164
$ws.Range("D82").Value = 166
$ws.Range("D83").Value = 167
$ws.Range("D84").Value = 165
$ws.Range("D85").Value = 163
$ws.Range("D88").Value = 61
$ws.Range("D90").Value = 62
$ws.Range("D91").Value = 63
$ws.Range("D92").Value = 64
$ws.Range("D93").Value = 65
$ws.Range("D94").Value = 56
$ws.Range("D95").Value = 55
$ws.Range("D96").Value = 66
$ws.Range("D97").Value = 58
$ws.Range("D98").Value = 59
$ws.Range("D99").Value = 60
$ws.Range("D102").Value = 130
$ws.Range("D103").Value = 131
$ws.Range("D104").Value = 115
$ws.Range("D105").Value = 132
$ws.Range("D106").Value = 133
$ws.Range("D107").Value = 119
$ws.Range("D108").Value = 114
$ws.Range("D109").Value = 135
$ws.Range("D110").Value = 120
$ws.Range("D111").Value = 117
$ws.Range("D112").Value = 122
$ws.Range("D113").Value = 126
$ws.Range("D114").Value = 124
$ws.Range("D115").Value = 127
$ws.Range("D116").Value = 116
$ws.Range("D117").Value = 128
$ws.Range("D118").Value = 129
$ws.Range("D119").Value = 123
$ws.Range("D121").Value = 121
$ws.Range("D122").Value = 118
$ws.Range("D123").Value = 125
$ws.Range("D125").Value = 191
$ws.Range("D126").Value = 178
$ws.Range("D127").Value = 192
$ws.Range("D128").Value = 193
$ws.Range("D129").Value = 194
$ws.Range("D130").Value = 195
$ws.Range("D131").Value = 196
$ws.Range("D132").Value = 182
$ws.Range("D133").Value = 181
$ws.Range("D134").Value = 183
$ws.Range("D135").Value = 185
$ws.Range("D136").Value = 187
$ws.Range("D137").Value = 188
$ws.Range("D138").Value = 180
$ws.Range("D139").Value = 189
$ws.Range("D140").Value = 179
$ws.Range("D141").Value = 190
$ws.Range("D142").Value = 186
$ws.Range("D143").Value = 184
$ws.Range("D144").Value = 177

$ws.Range("H45").Value = 5.1
$ws.Range("J45").Value = 6.779999999999999
